$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version property: 0.1.1 -> 0.2.0
$ws.Range("B3").Value = "0.2.0"

# Update Date property: 2023-10-20T07:19:33+00:00 -> 2023-10-20T08:59:58+00:00
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# Insert a new row 11 for the "Jurisdiction" property, pushing everything below down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (copy from the row now below, which still
# carries the original data-row style/border).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
